# Manager przenosi dwa kolejne elementy z listy "w_trakcie" do ListaElementow,
# oznaczajac je jako "Gotowe" (dopisanie wierszy 5 i 6).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = ""
$ws.Cells.Item(5, 3).Value = "prd.40062106.dld"
$ws.Cells.Item(5, 4).Value = "Gotowe"

$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = ""
$ws.Cells.Item(6, 3).Value = "prd.40662901siatka.dld"
$ws.Cells.Item(6, 4).Value = "Gotowe"
